$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows above row 297, shifting existing rows 297:358 down to 300:361
$xlShiftDown = -4121
$ws.Rows("297:299").Insert($xlShiftDown)

# New row 297 data
$ws.Range("A297").Value = 10
$ws.Range("B297").Value = "Vega Modelo de Temuco"
$ws.Range("C297").Value = "La Araucanía"
$ws.Range("D297").Value = 44951
$ws.Range("E297").Value = 9
$ws.Range("F297").Value = "Fruta"
$ws.Range("G297").Value = 100103
$ws.Range("H297").Value = "Frutos de hueso (carozo)"
$ws.Range("I297").Value = 100103004
$ws.Range("J297").Value = "Durazno"
$ws.Range("K297").Value = "Carson"
$ws.Range("L297").Value = "Primera"
$ws.Range("M297").Value = 150
$ws.Range("N297").Value = 25000
$ws.Range("O297").Value = 25000
$ws.Range("P297").Value = 25000
$ws.Range("Q297").Value = "$/bandeja 18 kilos granel"
$ws.Range("R297").Value = "Región de O'Higgins"
$ws.Range("S297").Value = 1389
$ws.Range("T297").Value = 18

# New row 298 data
$ws.Range("A298").Value = 10
$ws.Range("B298").Value = "Vega Modelo de Temuco"
$ws.Range("C298").Value = "La Araucanía"
$ws.Range("D298").Value = 44951
$ws.Range("E298").Value = 9
$ws.Range("F298").Value = "Fruta"
$ws.Range("G298").Value = 100103
$ws.Range("H298").Value = "Frutos de hueso (carozo)"
$ws.Range("I298").Value = 100103004
$ws.Range("J298").Value = "Durazno"
$ws.Range("K298").Value = "Carson"
$ws.Range("L298").Value = "Primera"
$ws.Range("M298").Value = 3
$ws.Range("N298").Value = 500000
$ws.Range("O298").Value = 500000
$ws.Range("P298").Value = 500000
$ws.Range("Q298").Value = "$/bins (400 kilos)"
$ws.Range("R298").Value = "Región de O'Higgins"
$ws.Range("S298").Value = 1250
$ws.Range("T298").Value = 400

# New row 299 data
$ws.Range("A299").Value = 10
$ws.Range("B299").Value = "Vega Modelo de Temuco"
$ws.Range("C299").Value = "La Araucanía"
$ws.Range("D299").Value = 44951
$ws.Range("E299").Value = 9
$ws.Range("F299").Value = "Fruta"
$ws.Range("G299").Value = 100103
$ws.Range("H299").Value = "Frutos de hueso (carozo)"
$ws.Range("I299").Value = 100103004
$ws.Range("J299").Value = "Durazno"
$ws.Range("K299").Value = "Springcrest"
$ws.Range("L299").Value = "Primera"
$ws.Range("M299").Value = 200
$ws.Range("N299").Value = 20000
$ws.Range("O299").Value = 20000
$ws.Range("P299").Value = 20000
$ws.Range("Q299").Value = "$/bandeja 18 kilos granel"
$ws.Range("R299").Value = "Región de O'Higgins"
$ws.Range("S299").Value = 1111
$ws.Range("T299").Value = 18

# Copy the date style (format) from the row below into the new rows' column D
$ws.Range("D300").Copy()
$ws.Range("D297:D299").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
